$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two description strings for Network Drive migration wording.
$ws.Range("B13").Value = "The target destination document library. Not used for Network Drive migrations (which uses the queue item meta data for this)"
$ws.Range("B14").Value = "The target destination SharePoint site. Not used for Network Drive migrations (which uses the queue item meta data for this)"

# Move the active selection from B20 to A4.
$ws.Range("A4").Select()
